$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'Douglas'
$ws.Range("B2").Value = 424
$ws.Range("C2").Value = 'Active'
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 7
$ws.Range("G2").Value = 3
$ws.Range("H2").Value = 10
$ws.Range("I2").Value = 27
$ws.Range("J2").Value = 5.4
$ws.Range("K2").Value = 18.07
$ws.Range("L2").Value = '04. Full House (3,2)'
$ws.Range("M2").Value = 'Ace,8'
$ws.Range("N2").Value = 18.07

# Row 3
$ws.Range("A3").Value = 'Alex'
$ws.Range("B3").Value = 348
$ws.Range("C3").Value = 'Active'
$ws.Range("D3").Value = 81
$ws.Range("E3").Value = 13
$ws.Range("F3").Value = 9
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 9
$ws.Range("I3").Value = 274
$ws.Range("J3").Value = 3.38
$ws.Range("K3").Value = 8.07
$ws.Range("L3").Value = '03. 4 of a Kind'
$ws.Range("M3").Value = 'King'
$ws.Range("N3").Value = 11.02

# Row 4
$ws.Range("A4").Value = 'Andy'
$ws.Range("B4").Value = 349
$ws.Range("C4").Value = 'Active'
$ws.Range("D4").Value = 191
$ws.Range("E4").Value = 19
$ws.Range("F4").Value = 11
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 11
$ws.Range("I4").Value = 790
$ws.Range("J4").Value = 4.14
$ws.Range("K4").Value = 11.04
$ws.Range("L4").Value = '03. 4 of a Kind'
$ws.Range("M4").Value = 'King'
$ws.Range("N4").Value = 16.05

# Row 5
$ws.Range("A5").Value = 'Anthony'
$ws.Range("B5").Value = 350
$ws.Range("C5").Value = 'Active'
$ws.Range("D5").Value = 114
$ws.Range("E5").Value = 13
$ws.Range("F5").Value = 10
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 10
$ws.Range("I5").Value = 453
$ws.Range("J5").Value = 3.97
$ws.Range("K5").Value = 9.08
$ws.Range("L5").Value = '02. Straight Flush'
$ws.Range("M5").Value = '7,8,9,10,Jack,(H)'
$ws.Range("N5").Value = 12.12

# Row 6
$ws.Range("A6").Value = 'Crafty'
$ws.Range("B6").Value = 354
$ws.Range("C6").Value = 'Active'
$ws.Range("D6").Value = 17
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 11
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 11
$ws.Range("I6").Value = 56
$ws.Range("J6").Value = 3.29
$ws.Range("K6").Value = 14.01
$ws.Range("L6").Value = '03. 4 of a Kind'
$ws.Range("M6").Value = 2
$ws.Range("N6").Value = 14.03

# Row 7
$ws.Range("A7").Value = 'Illya'
$ws.Range("B7").Value = 355
$ws.Range("C7").Value = 'Active'
$ws.Range("D7").Value = 28
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 10
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 10
$ws.Range("I7").Value = 111
$ws.Range("J7").Value = 3.96
$ws.Range("K7").Value = 4.09
$ws.Range("L7").Value = '05. Flush'
$ws.Range("M7").Value = 'Ace,4,8,9,King,(H)'
$ws.Range("N7").Value = 17.07

# Row 8
$ws.Range("A8").Value = 'Jon'
$ws.Range("B8").Value = 357
$ws.Range("C8").Value = 'Active'
$ws.Range("D8").Value = 187
$ws.Range("E8").Value = 19
$ws.Range("F8").Value = 13
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 13
$ws.Range("I8").Value = 682
$ws.Range("J8").Value = 3.65
$ws.Range("K8").Value = 13.08
$ws.Range("L8").Value = '03. 4 of a Kind'
$ws.Range("M8").Value = 'Ace'
$ws.Range("N8").Value = 12.04

# Row 9
$ws.Range("A9").Value = 'Maisy'
$ws.Range("B9").Value = 360
$ws.Range("C9").Value = 'Active'
$ws.Range("D9").Value = 115
$ws.Range("E9").Value = 13
$ws.Range("F9").Value = 9
$ws.Range("G9").Value = 3
$ws.Range("H9").Value = 12
$ws.Range("I9").Value = 431
$ws.Range("J9").Value = 3.75
$ws.Range("K9").Value = 18.04
$ws.Range("L9").Value = '02. Straight Flush'
$ws.Range("M9").Value = '6,7,8,9,10,(D)'
$ws.Range("N9").Value = 15.09

# Row 10
$ws.Range("A10").Value = 'Mark'
$ws.Range("B10").Value = 361
$ws.Range("C10").Value = 'Active'
$ws.Range("D10").Value = 134
$ws.Range("E10").Value = 15
$ws.Range("F10").Value = 9
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 9
$ws.Range("I10").Value = 475
$ws.Range("J10").Value = 3.54
$ws.Range("K10").Value = 14.12
$ws.Range("L10").Value = '03. 4 of a Kind'
$ws.Range("M10").Value = 'Ace'
$ws.Range("N10").Value = 14.02

# Row 11
$ws.Range("A11").Value = 'Matt'
$ws.Range("B11").Value = 362
$ws.Range("C11").Value = 'Active'
$ws.Range("D11").Value = 182
$ws.Range("E11").Value = 19
$ws.Range("F11").Value = 11
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 11
$ws.Range("I11").Value = 693
$ws.Range("J11").Value = 3.81
$ws.Range("K11").Value = 13.12
$ws.Range("L11").Value = '03. 4 of a Kind'
$ws.Range("M11").Value = 'King'
$ws.Range("N11").Value = 16.02

# Row 12
$ws.Range("A12").Value = 'Pepe'
$ws.Range("B12").Value = 364
$ws.Range("C12").Value = 'Active'
$ws.Range("D12").Value = 97
$ws.Range("E12").Value = 12
$ws.Range("F12").Value = 9
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 10
$ws.Range("I12").Value = 358
$ws.Range("J12").Value = 3.69
$ws.Range("K12").Value = 17.02
$ws.Range("L12").Value = '03. 4 of a Kind'
$ws.Range("M12").Value = 'Ace'
$ws.Range("N12").Value = 15.08

# Row 13
$ws.Range("A13").Value = 'Prashant'
$ws.Range("B13").Value = 365
$ws.Range("C13").Value = 'Active'
$ws.Range("D13").Value = 30
$ws.Range("E13").Value = 5
$ws.Range("F13").Value = 9
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 9
$ws.Range("I13").Value = 117
$ws.Range("J13").Value = 3.9
$ws.Range("K13").Value = 16.01
$ws.Range("L13").Value = '03. 4 of a Kind'
$ws.Range("M13").Value = 'King'
$ws.Range("N13").Value = 16.01

# Row 14
$ws.Range("A14").Value = 'Richard'
$ws.Range("B14").Value = 366
$ws.Range("C14").Value = 'Active'
$ws.Range("D14").Value = 134
$ws.Range("E14").Value = 19
$ws.Range("F14").Value = 12
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 12
$ws.Range("I14").Value = 561
$ws.Range("J14").Value = 4.19
$ws.Range("K14").Value = 12.05
$ws.Range("L14").Value = '03. 4 of a Kind'
$ws.Range("M14").Value = 5
$ws.Range("N14").Value = 18.05

# Row 15
$ws.Range("A15").Value = 'Richard Snr'
$ws.Range("B15").Value = 367
$ws.Range("C15").Value = 'Inactive'
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 1
$ws.Range("K15").Value = 3.11
$ws.Range("L15").ClearContents()
$ws.Range("M15").ClearContents()
$ws.Range("N15").ClearContents()

# Row 16
$ws.Range("A16").Value = 'Stuart'
$ws.Range("B16").Value = 368
$ws.Range("C16").Value = 'Inactive'
$ws.Range("D16").Value = 48
$ws.Range("E16").Value = 7
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 12
$ws.Range("I16").Value = 183
$ws.Range("J16").Value = 3.81
$ws.Range("K16").Value = 4.02
$ws.Range("L16").ClearContents()
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()

# Row 17
$ws.Range("A17").Value = 'Matthew'
$ws.Range("B17").Value = 363
$ws.Range("C17").Value = 'Inactive'
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 3
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 3
$ws.Range("I17").Value = 5
$ws.Range("J17").Value = 2.5
$ws.Range("K17").Value = 4.06
$ws.Range("L17").ClearContents()
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()

# Row 18
$ws.Range("A18").Value = 'Jonathan'
$ws.Range("B18").Value = 358
$ws.Range("C18").Value = 'Inactive'
$ws.Range("D18").Value = 10
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 7
$ws.Range("I18").Value = 35
$ws.Range("J18").Value = 3.5
$ws.Range("K18").Value = 4.07
$ws.Range("L18").ClearContents()
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()

# Row 19
$ws.Range("A19").Value = 'Keith'
$ws.Range("B19").Value = 359
$ws.Range("C19").Value = 'Inactive'
$ws.Range("D19").Value = 98
$ws.Range("E19").Value = 12
$ws.Range("F19").Value = 10
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 10
$ws.Range("I19").Value = 333
$ws.Range("J19").Value = 3.4
$ws.Range("K19").Value = 13.04
$ws.Range("L19").Value = '03. 4 of a Kind'
$ws.Range("M19").Value = 8
$ws.Range("N19").Value = 12.06

# Row 20
$ws.Range("A20").Value = 'Jim'
$ws.Range("B20").Value = 356
$ws.Range("C20").Value = 'Inactive'
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 9
$ws.Range("I20").Value = 17
$ws.Range("J20").Value = 5.67
$ws.Range("K20").Value = 7.05
$ws.Range("L20").ClearContents()
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()

# Row 21
$ws.Range("A21").Value = 'Bob'
$ws.Range("B21").Value = 351
$ws.Range("C21").Value = 'Inactive'
$ws.Range("D21").Value = 35
$ws.Range("E21").Value = 9
$ws.Range("F21").Value = 11
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 11
$ws.Range("I21").Value = 139
$ws.Range("J21").Value = 3.97
$ws.Range("K21").Value = 7.01
$ws.Range("L21").Value = '04. Full House (3,2)'
$ws.Range("M21").Value = '4,5'
$ws.Range("N21").Value = 10.12

# Row 22
$ws.Range("A22").Value = 'Chris'
$ws.Range("B22").Value = 353
$ws.Range("C22").Value = 'Inactive'
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 3
$ws.Range("I22").Value = 3
$ws.Range("J22").Value = 3
$ws.Range("K22").Value = 5.05
$ws.Range("L22").ClearContents()
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

# Row 23
$ws.Range("A23").Value = 'John'
$ws.Range("B23").Value = 113
$ws.Range("C23").Value = 'Inactive'
$ws.Range("D23").Value = 9
$ws.Range("E23").Value = 3
$ws.Range("F23").Value = 12
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 12
$ws.Range("I23").Value = 46
$ws.Range("J23").Value = 5.11
$ws.Range("K23").Value = 5.06
$ws.Range("L23").ClearContents()
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()
